# Updated cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.060.79"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.03"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.06"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4221"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07204"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8416"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.78"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.836.30"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.668"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07071"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.284"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.62"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008754"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.90"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.130.09"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.047.34"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.994"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.95"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.248"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.23"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.286"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.52"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08736"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.182"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7394"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.946"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.421"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9998"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.089"
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05247"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.333"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.867"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1688"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5049"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.602"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.55"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.41"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4718"
$ws.Range("E47").Value = "  -1.69%  "

# Rows 48 and 49 swapped (RenderToken moved up to rank 46, PaxDollar moved down to rank 47)
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.934"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9997"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06340"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.651"
$ws.Range("E51").Value = "  -2.53%  "
